$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell address -> new text value (dots replaced by commas in percentages)
$updates = @{
    "C2" = "55,8%"
    "E2" = "44,2%"
    "G2" = "9,3%"
    "C3" = "52,4%"
    "E3" = "47,6%"
    "G3" = "68,9%"
    "C4" = "45,8%"
    "E4" = "54,2%"
    "G4" = "12,3%"
    "C5" = "41,3%"
    "E5" = "58,7%"
    "G5" = "3,9%"
    "C6" = "38,1%"
    "E6" = "61,9%"
    "G6" = "2,3%"
    "G7" = "1,2%"
    "C8" = "39,9%"
    "E8" = "60,1%"
    "G8" = "2,1%"
    "C9" = "50,7%"
    "E9" = "49,3%"
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
